# Updates crypto price/volume table cells per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'27.206.74"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.75%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'1.569.46"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.51%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.53%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'211.66"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.01%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +0.65%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.40%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.55%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +0.37%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +0.60%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0868"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.02%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'1.792.92"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.45%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'1.585.20"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.47%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +0.75%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.520"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.31%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'27.208.02"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.75%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'62.33"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.42%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -0.47%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "BitcoinCash"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(19, 4).Value = "'216.33"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.37%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "Chainlink"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(20, 4).Value = "'7.43"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.75%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +0.44%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +1.18%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -0.03%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'1.95"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.06%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'154.03"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.55%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'6.67"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.50%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'15.10"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.37%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'0.106"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.89%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.48%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +2.54%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.34%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Maker"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(33, 4).Value = "'1.453.48"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +2.13%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(34, 4).Value = "'3.17"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +2.09%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +5.18%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +0.24%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +1.38%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +0.98%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +0.81%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -0.02%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +0.48%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.69%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -0.15%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'64.74"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.32%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'1.74"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.69%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'1.707.06"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.36%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'85.89"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -1.98%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +3.13%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -0.19%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'0.0959"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.20%  "
